$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: editUserExp1 / reuse existing xpath prefix string
$ws.Range("A16").Value = "editUserExp1"
$ws.Range("B16").Value = "//tr[td[@class='sorting_1' ]/label[contains(text(),'"

# Row 17: editUserExp2 / new xpath suffix string (needs leading apostrophe preserved,
# so an extra leading apostrophe is supplied - Excel consumes the first as a
# quote-prefix marker and keeps the rest, applying the quote-prefix cell style)
$ws.Range("A17").Value = "editUserExp2"
$ws.Range("B17").Value = "'') ]]/td/button[@id='edit']"

# Row 18: editUser / user_save_btn (set B18 before A18 so the shared-string
# table ends up with "user_save_btn" allocated before "editUser")
$ws.Range("B18").Value = "user_save_btn"
$ws.Range("A18").Value = "editUser"

# Leave the active selection on A18, matching the saved workbook state
$null = $ws.Range("A18").Select()
